$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "10×73=730"
$t.Cell(1,2).Range.Text = "16×92=1472"
$t.Cell(1,3).Range.Text = "48×13=624"
$t.Cell(1,4).Range.Text = "73×48=3504"
$t.Cell(1,5).Range.Text = "22×59=1298"
$t.Cell(2,1).Range.Text = "61×86=5246"
$t.Cell(2,2).Range.Text = "90×51=4590"
$t.Cell(2,3).Range.Text = "85×65=5525"
$t.Cell(2,4).Range.Text = "33×54=1782"
$t.Cell(2,5).Range.Text = "45×44=1980"
$t.Cell(3,1).Range.Text = "29×84=2436"
$t.Cell(3,2).Range.Text = "76×56=4256"
$t.Cell(3,3).Range.Text = "100×49=4900"
$t.Cell(3,4).Range.Text = "79×96=7584"
$t.Cell(3,5).Range.Text = "25×19=475"
$t.Cell(4,1).Range.Text = "96×91=8736"
$t.Cell(4,2).Range.Text = "71×37=2627"
$t.Cell(4,3).Range.Text = "48×53=2544"
$t.Cell(4,4).Range.Text = "25×54=1350"
$t.Cell(4,5).Range.Text = "26×84=2184"
$t.Cell(5,1).Range.Text = "34×79=2686"
$t.Cell(5,2).Range.Text = "24×44=1056"
$t.Cell(5,3).Range.Text = "95×51=4845"
$t.Cell(5,4).Range.Text = "56×84=4704"
$t.Cell(5,5).Range.Text = "24×51=1224"
$t.Cell(6,1).Range.Text = "93×36=3348"
$t.Cell(6,2).Range.Text = "85×72=6120"
$t.Cell(6,3).Range.Text = "24×75=1800"
$t.Cell(6,4).Range.Text = "83×23=1909"
$t.Cell(6,5).Range.Text = "77×91=7007"
$t.Cell(7,1).Range.Text = "39×52=2028"
$t.Cell(7,2).Range.Text = "43×20=860"
$t.Cell(7,3).Range.Text = "38×60=2280"
$t.Cell(7,4).Range.Text = "37×83=3071"
$t.Cell(7,5).Range.Text = "44×15=660"
$t.Cell(8,1).Range.Text = "86×13=1118"
$t.Cell(8,2).Range.Text = "54×30=1620"
$t.Cell(8,3).Range.Text = "95×14=1330"
$t.Cell(8,4).Range.Text = "26×53=1378"
$t.Cell(8,5).Range.Text = "80×31=2480"
$t.Cell(9,1).Range.Text = "94×88=8272"
$t.Cell(9,2).Range.Text = "93×76=7068"
$t.Cell(9,3).Range.Text = "72×36=2592"
$t.Cell(9,4).Range.Text = "100×100=10000"
$t.Cell(9,5).Range.Text = "81×29=2349"
$t.Cell(10,1).Range.Text = "78×83=6474"
$t.Cell(10,2).Range.Text = "29×87=2523"
$t.Cell(10,3).Range.Text = "63×83=5229"
$t.Cell(10,4).Range.Text = "48×88=4224"
$t.Cell(10,5).Range.Text = "64×31=1984"
$t.Cell(11,1).Range.Text = "55×67=3685"
$t.Cell(11,2).Range.Text = "88×19=1672"
$t.Cell(11,3).Range.Text = "84×20=1680"
$t.Cell(11,4).Range.Text = "24×62=1488"
$t.Cell(11,5).Range.Text = "53×44=2332"
$t.Cell(12,1).Range.Text = "29×31=899"
$t.Cell(12,2).Range.Text = "23×48=1104"
$t.Cell(12,3).Range.Text = "34×76=2584"
$t.Cell(12,4).Range.Text = "65×78=5070"
$t.Cell(12,5).Range.Text = "69×82=5658"
$t.Cell(13,1).Range.Text = "22×85=1870"
$t.Cell(13,2).Range.Text = "31×48=1488"
$t.Cell(13,3).Range.Text = "59×34=2006"
$t.Cell(13,4).Range.Text = "46×87=4002"
$t.Cell(13,5).Range.Text = "86×20=1720"
$t.Cell(14,1).Range.Text = "60×85=5100"
$t.Cell(14,2).Range.Text = "88×50=4400"
$t.Cell(14,3).Range.Text = "10×18=180"
$t.Cell(14,4).Range.Text = "73×70=5110"
$t.Cell(14,5).Range.Text = "57×37=2109"
$t.Cell(15,1).Range.Text = "79×78=6162"
$t.Cell(15,2).Range.Text = "77×85=6545"
$t.Cell(15,3).Range.Text = "59×48=2832"
$t.Cell(15,4).Range.Text = "55×38=2090"
$t.Cell(15,5).Range.Text = "75×90=6750"
$t.Cell(16,1).Range.Text = "38×43=1634"
$t.Cell(16,2).Range.Text = "71×55=3905"
$t.Cell(16,3).Range.Text = "85×57=4845"
$t.Cell(16,4).Range.Text = "82×63=5166"
$t.Cell(16,5).Range.Text = "30×38=1140"
$t.Cell(17,1).Range.Text = "31×23=713"
$t.Cell(17,2).Range.Text = "24×39=936"
$t.Cell(17,3).Range.Text = "53×75=3975"
$t.Cell(17,4).Range.Text = "81×18=1458"
$t.Cell(17,5).Range.Text = "32×46=1472"
$t.Cell(18,1).Range.Text = "36×20=720"
$t.Cell(18,2).Range.Text = "96×78=7488"
$t.Cell(18,3).Range.Text = "28×32=896"
$t.Cell(18,4).Range.Text = "27×86=2322"
$t.Cell(18,5).Range.Text = "91×41=3731"
$t.Cell(19,1).Range.Text = "14×41=574"
$t.Cell(19,2).Range.Text = "50×29=1450"
$t.Cell(19,3).Range.Text = "12×59=708"
$t.Cell(19,4).Range.Text = "95×25=2375"
$t.Cell(19,5).Range.Text = "97×24=2328"
$t.Cell(20,1).Range.Text = "80×11=880"
$t.Cell(20,2).Range.Text = "81×89=7209"
$t.Cell(20,3).Range.Text = "56×11=616"
$t.Cell(20,4).Range.Text = "22×57=1254"
$t.Cell(20,5).Range.Text = "57×52=2964"
